$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("answers")

$ws.Range("B2").Value = "Am I the only question left?"
$ws.Range("B3").Value = "No I'm here too"

$ws.Range("A2").Value = 1
$ws.Range("C2").Value = $false
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

$ws.Range("A3").Value = 1
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

$ws.Range("D4").Select()
